$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet gets a new blank column inserted before column N
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.17

# Move the selection on the Repayment schedule sheet
$ws.Range("S7").Select()

# Make "Repayment schedule" the active sheet/tab (this also moves tabSelected
# off the previously active "Summary" sheet and bumps workbook activeTab)
$ws.Activate()
